$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the paragraph that currently contains *only* a manual page break
# (the form-feed character) right after the "Consult with Lena..." bullet,
# and just before the "Progress" Heading1 paragraph.
# ---------------------------------------------------------------------------
$breakIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.Length -eq 2 -and [int][char]$t[0] -eq 12) {
        $breakIdx = $i
        break
    }
}

if ($breakIdx -eq -1) {
    throw "Could not find the lone page-break paragraph"
}

# ---------------------------------------------------------------------------
# Insert a brand-new, empty paragraph right before it — this becomes the new
# "TODO" Heading1 paragraph. After the insert, the new empty paragraph sits
# at $breakIdx and the original page-break paragraph moves to $breakIdx + 1.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item($breakIdx).Range.InsertParagraphBefore()

$todoPara = $d.Paragraphs.Item($breakIdx)
$todoPara.Range.Text = "TODO"
$todoPara.Style = "Heading 1"

$figIdx = $breakIdx + 1

# ---------------------------------------------------------------------------
# Turn the old break-only paragraph into the new bulleted paragraph:
#   "Figure out how to combine all species into a single shapefile – 1
#    each for each season and 1 for year round" + page break
# Re-use the same list (numId 1 / ilvl 0) as the other top-level bullets.
# ---------------------------------------------------------------------------
$figPara = $d.Paragraphs.Item($figIdx)
$figPara.Style = "List Paragraph"

$listTemplate = $d.Paragraphs.Item(7).Range.ListFormat.ListTemplate
$figPara.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)

$figPara = $d.Paragraphs.Item($figIdx)
$insertPoint = $d.Range($figPara.Range.Start, $figPara.Range.Start)
$insertPoint.InsertBefore("Figure out how to combine all species into a single shapefile " + [char]0x2013 + " 1 each for each season and 1 for year round")

Write-Output "done"
